$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; temporarily unprotect to make edits, then restore protection.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A42).
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.
Model holdings provided as of 2021-05-18 for illustrative purposes only and are subject to change."

# Update model holdings Weight (D) and Percent Change (E) values.
$ws.Range("D2").Value = 0.05759297531619011
$ws.Range("E2").Value = -0.01124574324859418
$ws.Range("D3").Value = 0.05201348048211011
$ws.Range("E3").Value = -0.008565135818582315
$ws.Range("D4").Value = 0.3139307790905724
$ws.Range("E4").Value = 0.003311258278145601
$ws.Range("D5").Value = 0.03468969052000327
$ws.Range("E5").Value = -0.01165304443812498
$ws.Range("D6").Value = 0.03122463552589876
$ws.Range("E6").Value = -0.003842077371489117
$ws.Range("D7").Value = 0.031440386034302
$ws.Range("E7").Value = -0.01408878362786181
$ws.Range("D8").Value = 0.02891778102405331
$ws.Range("E8").Value = 0.0003521333411584671
$ws.Range("D9").Value = 0.02357175072733591
$ws.Range("E9").Value = 0.02174382604939162
$ws.Range("D10").Value = 0.024279039021354
$ws.Range("E10").Value = -0.01155566817538412
$ws.Range("D11").Value = 0.02342303774170103
$ws.Range("E11").Value = -0.01743485703417236
$ws.Range("D12").Value = 0.02357429645604863
$ws.Range("E12").Value = -0.01357042583060375
$ws.Range("D13").Value = 0.01984480389191006
$ws.Range("E13").Value = -0.00235183443085607
$ws.Range("D14").Value = 0.02191395097520503
$ws.Range("E14").Value = -0.01023742104116743
$ws.Range("D15").Value = 0.02036646613395879
$ws.Range("E15").Value = -0.01018718165057342
$ws.Range("D16").Value = 0.02167210674749638
$ws.Range("E16").Value = 0.004282602843648231
$ws.Range("D17").Value = 0.019267984194419
$ws.Range("E17").Value = -0.0006330856041839583
$ws.Range("D18").Value = 0.01408551696749401
$ws.Range("E18").Value = -0.008855955178022734
$ws.Range("D19").Value = 0.01647022833913689
$ws.Range("E19").Value = -0.009106482733748034
$ws.Range("D20").Value = 0.01536453350157772
$ws.Range("E20").Value = -0.01311701760441841
$ws.Range("D21").Value = 0.01715121076979018
$ws.Range("E21").Value = -0.02830036983437856
$ws.Range("D22").Value = 0.01223710577799803
$ws.Range("E22").Value = 0.001802957543816985
$ws.Range("D23").Value = 0.01506901682684251
$ws.Range("E23").Value = -0.005490483162518323
$ws.Range("D24").Value = 0.01430816215782754
$ws.Range("E24").Value = -0.05801721389862924
$ws.Range("D25").Value = 0.01403863313036804
$ws.Range("E25").Value = -0.0001888930865130423
$ws.Range("D26").Value = 0.01362113362148154
$ws.Range("E26").Value = 0.002055850608189136
$ws.Range("D27").Value = 0.01302235701384325
$ws.Range("E27").Value = -0.006231214720328504
$ws.Range("D28").Value = 0.01410460993283943
$ws.Range("E28").Value = -0.04178323256023819
$ws.Range("D29").Value = 0.01446546697786786
$ws.Range("E29").Value = -0.001495886312640415
$ws.Range("D30").Value = 0.01342108177347342
$ws.Range("E30").Value = -0.004931714719271674
$ws.Range("D31").Value = 0.01240363886462196
$ws.Range("E31").Value = -0.004515290415269879
$ws.Range("D32").Value = 0.01363863550638151
$ws.Range("E32").Value = 0.002737616562580092
$ws.Range("D33").Value = 0.01270795951782201
$ws.Range("E33").Value = -0.006134969325153672
$ws.Range("D34").Value = 0.00601025334667861
$ws.Range("E34").Value = -0.01057145882602095
$ws.Range("D35").Value = 0.005186285819994069
$ws.Range("E35").Value = -0.005440340328056625
$ws.Range("D36").Value = 0.005183952235340739
$ws.Range("E36").Value = -0.004706171222786049
$ws.Range("D37").Value = 0.005120521161582069
$ws.Range("E37").Value = -0.004391597961635685
$ws.Range("D38").Value = 0.004666532874479876
$ws.Range("E38").Value = -0.01409283084056911
$ws.Range("D39").Value = 0.9999999999999998
$ws.Range("E39").Value = -0.004791273581403299

# Restore sheet protection.
$ws.Protect()
